# Introduce a simple raycast based groundcheck
# - Add a new "Sheet2" after the existing "Sheet1"
# - Populate it with a header row (re-using existing weekday strings) and
#   three new task notes
# - Make it the active sheet/tab
# - Update Sheet1's view state (scroll/selection) to no longer be the
#   selected tab

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Create Sheet2 right after Sheet1
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# ---------------------------------------------------------------------
# 2. Header row - reuses the Monday..Friday strings already in the
#    workbook's shared string table
# ---------------------------------------------------------------------
$ws2.Range("B1").Value = "Monday"
$ws2.Range("C1").Value = "Tuesday"
$ws2.Range("D1").Value = "Wednesday"
$ws2.Range("E1").Value = "Thursday"
$ws2.Range("F1").Value = "Friday"
$ws2.Range("G1").Value = "Monday"
$ws2.Range("H1").Value = "Tuesday"
$ws2.Range("I1").Value = "Wednesday"
$ws2.Range("J1").Value = "Thursday"
$ws2.Range("K1").Value = "Friday"
$ws2.Range("K1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 3. New task notes (new shared strings, appended to the table)
# ---------------------------------------------------------------------
$ws2.Range("B2").Value = "Implement a raycast based groundcheck"
$ws2.Range("C2").Value = "Implement a raycast based groundcheck"
$ws2.Range("C3").Value = "Start designing levels"
$ws2.Range("C4").Value = "Credit Controller Code creators"

# ---------------------------------------------------------------------
# 4. Column widths
# ---------------------------------------------------------------------
$ws2.Columns.Item(2).ColumnWidth = 37.140625
$ws2.Columns.Item(3).ColumnWidth = 37.140625
$ws2.Columns.Item(4).ColumnWidth = 26.5703125
$ws2.Columns.Item(5).ColumnWidth = 29.28515625
$ws2.Columns.Item(6).ColumnWidth = 20.85546875
$ws2.Columns.Item(7).ColumnWidth = 25.140625
$ws2.Columns.Item(8).ColumnWidth = 27.85546875
$ws2.Columns.Item(9).ColumnWidth = 28.5703125
$ws2.Columns.Item(10).ColumnWidth = 20.28515625
$ws2.Columns.Item(11).ColumnWidth = 20

# ---------------------------------------------------------------------
# 5. View state: Sheet2 becomes the active tab / selected cell C5
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("C5").Select()

# ---------------------------------------------------------------------
# 6. Sheet1 view state: scroll / selection moves, no longer tabSelected
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("F1").Select()
$ws1.Range("J3").Select()

# Leave Sheet2 as the active sheet/tab
$ws2.Activate()
